# Convert "Manufacturing Employees on Avg." columns (AF:AO, years 2008-2017)
# from a cumulative total over all manufacturing HUBZone businesses in the
# county/year to a true average, by dividing by the corresponding count of
# HUBZone businesses for that year (columns B:K, years 2008-2017).
#
# Column mapping (same year): B<->AF, C<->AG, D<->AH, E<->AI, F<->AJ,
#                              G<->AK, H<->AL, I<->AM, J<->AN, K<->AO

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$maxRow = $usedRange.Rows.Count

$countCols = 2, 3, 4, 5, 6, 7, 8, 9, 10, 11          # B .. K  (year count of HUBZone businesses)
$empCols   = 32, 33, 34, 35, 36, 37, 38, 39, 40, 41  # AF .. AO (Manufacturing Employees on Avg.)

# Row 1 is the header row, data starts at row 2.
for ($r = 2; $r -le $maxRow; $r++) {
    for ($i = 0; $i -lt $countCols.Length; $i++) {
        $countCol = $countCols[$i]
        $empCol = $empCols[$i]

        $countCell = $ws.Cells.Item($r, $countCol)
        $empCell = $ws.Cells.Item($r, $empCol)

        $countVal = $countCell.Value2
        $empVal = $empCell.Value2

        if ($countVal -ne $null -and $empVal -ne $null) {
            if ($countVal -gt 1 -and $empVal -ne 0) {
                $empCell.Value2 = $empVal / $countVal
            }
        }
    }
}
